$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New sound-duration rows (Sword_Draw_01..05) filling previously-empty rows 16-20
$names = @("Sword_Draw_01", "Sword_Draw_02", "Sword_Draw_03", "Sword_Draw_04", "Sword_Draw_05")
$durations = @(1.39, 1.34, 1.25, 1.2, 1.22)

for ($i = 0; $i -lt 5; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $durations[$i]
    $ws.Cells.Item($row, 4).Value = 20
}

# Update the view: scroll so row 8 is the top-left visible row, and move the
# active selection to B20 (matching the post-edit sheetView state)
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("B20").Select()
